$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows: row, dateSerial, B, C, D
$data = @(
    @(302, 44376, 0, 5, 58.91363261458702),
    @(303, 44377, 0, 5, 58.91363261458702),
    @(304, 44378, 0, 1, 11.7827265229174),
    @(305, 44379, 0, 1, 11.7827265229174),
    @(306, 44380, 0, 0, 0),
    @(307, 44381, 0, 0, 0),
    @(308, 44382, 0, 0, 0),
    @(309, 44383, 0, 0, 0),
    @(310, 44384, 0, 0, 0),
    @(311, 44385, 0, 0, 0),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 0, 0, 0),
    @(317, 44391, 0, 0, 0),
    @(318, 44392, 0, 0, 0),
    @(319, 44393, 0, 0, 0),
    @(320, 44394, 0, 0, 0),
    @(321, 44395, 0, 0, 0),
    @(322, 44396, 0, 0, 0),
    @(323, 44397, 0, 0, 0),
    @(324, 44398, 0, 0, 0),
    @(325, 44399, 0, 0, 0),
    @(326, 44400, 1, 1, 11.7827265229174),
    @(327, 44401, 0, 1, 11.7827265229174),
    @(328, 44402, 0, 1, 11.7827265229174)
)

foreach ($item in $data) {
    $r = $item[0]
    $prevR = $r - 1
    # Copy formatting (style) from the row above down to the new row
    $ws.Range("A" + $prevR + ":D" + $prevR).Copy($ws.Range("A" + $r + ":D" + $r))
    $ws.Range("A" + $r).Value = $item[1]
    $ws.Range("B" + $r).Value = $item[2]
    $ws.Range("C" + $r).Value = $item[3]
    $ws.Range("D" + $r).Value = $item[4]
}
